# Generate Report for Handoff
#
# A new handoff was generated for file "b37ee5e8-b13a-487a-a73a-599a8dda12da.md"
# (row 7 in each of the report tables), refreshing its handoff timestamps:
#   - zh-cn handoff happened at  2016-08-13 16:56:55
#   - de-de handoff happened at  2016-08-13 16:57:05
# and the Overview sheet's "Latest HO Xliff Generate Date" (which mirrors the
# de-de "Latest Handoff Datetime") is refreshed to match.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G7 -> Latest HO Xliff Generate Date for b37ee5e8-...
$wsOverview.Range("G7").Value = "2016-08-13 16:57:05"

# zh-cn!H7 -> Latest Handoff Datetime for b37ee5e8-...
$wsZhCn.Range("H7").Value = "2016-08-13 16:56:55"

# de-de!H7 -> Latest Handoff Datetime for b37ee5e8-...
$wsDeDe.Range("H7").Value = "2016-08-13 16:57:05"
